$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Cashflow")

# Insert a new column at E (shifts Novelty's neighbours -> Values0.. right by one)
$ws.Columns("E:E").Insert()

# Header for the newly inserted column
$ws.Range("E1").Value = "AccidentYear"

# Keep the new column's width close to its neighbours (9 chars), but not "best fit"
$ws.Columns("E:E").ColumnWidth = 8.14

# Grow the table to cover the new column
$lo = $ws.ListObjects.Item("Table_Cashflow")
$lastCol = $lo.ListColumns.Count + 1
$lo.Resize($ws.Range("A1:BN9"))

# Re-sync every table column name from its (now correctly shifted) header cell so the
# table definition lines back up with the worksheet headers column-for-column.
for ($col = 6; $col -le 66; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $cell.Value = $cell.Text
}

# The Novelty validation list now also covers the new AccidentYear column
$ws.Range("D2:E9").Validation.Delete()
$ws.Range("D2:E9").Validation.Add(3, 1, 1, "Novelty_SystemName")
$ws.Range("D2:E9").Validation.IgnoreBlank = $true
$ws.Range("D2:E9").Validation.InCellDropdown = $true
$ws.Range("D2:E9").Validation.ShowInput = $false
$ws.Range("D2:E9").Validation.ShowError = $false

# Reset the view: normal zoom, selection on the new column's first data cell
$ws.Activate()
$excel.ActiveWindow.Zoom = 100
$ws.Range("E2").Select()
